# Carrera-Mapping.xlsx edit script
# - rename sheets
# - restructure "test"/"Vorlage zur Erstellung CSV" sheet: add bold header row,
#   replace the free-text "section_1"/"section_2" names with a numeric Barcode
#   column, and protect the "Tabelle1"/"Schienentypen" sheet.

$wb = $excel.ActiveWorkbook

# --- rename the worksheets -------------------------------------------------
$wsTypes = $wb.Worksheets.Item(1)
$wsTypes.Name = "Schienentypen"

$wsCsv = $wb.Worksheets.Item(2)
$wsCsv.Name = "Vorlage zur Erstellung CSV"

# --- rebuild the data sheet -------------------------------------------------
# Existing layout (rows 2-16):
#   2: Type | name | anzahl | pos_x | pos_y | rot_z           (header)
#   3: SECTION_START | section_1 | 7 | 5 | 0 | 0
#   4-6: straight_standard
#   7-9: curve_R1_30_outer
#   10: straight_standard
#   11: SECTION_END
#   12: SECTION_START | section_2 | 3 | 580 | 720 | 30
#   13: straight_standard
#   14-15: straight_standard_third
#   16: SECTION_END
#
# New layout (rows 1-15): same data, shifted up one row, with a new bold
# header row inserted at the top and the section "name" replaced by a
# sequential "Barcode" number.

$wsCsv.Cells.Item(1, 1).EntireRow.Delete() | Out-Null

$headerValues = @("Type", "Barcode", "anzahl", "pos_x", "pos_y", "rot_z")
for ($c = 0; $c -lt $headerValues.Length; $c++) {
    $wsCsv.Cells.Item(1, $c + 1).Value = $headerValues[$c]
}
$headerRange = $wsCsv.Range("A1:F1")
$headerRange.Font.Bold = $true

$wsCsv.Cells.Item(2, 2).Value = 1
$wsCsv.Cells.Item(11, 2).Value = 2

# --- protect the Schienentypen sheet ---------------------------------------
$wsTypes.Range("F18").Select()
$wsTypes.Protect()

# leave the CSV template sheet selected/active, matching the saved file
$wsCsv.Activate()
$wsCsv.Range("J11").Select()
